# Auto-generated edit script: updates cryptos list price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.772.63"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.582.92"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'585.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'169.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "2.583.07"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'26.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").Value = "66.588.22"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "2.551.76"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "'11.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.07%  "
$ws.Range("D20").Value = "'7.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("D21").Value = "'351.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "'4.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'69.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").Value = "'9.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.95%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "0.0₃0994"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").Value = "'533.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").Value = "'8.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").Value = "'0.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "'156.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").Value = "'18.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0288"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'149.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").Value = "'0.569"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  -0.89%  "
